$wb = $excel.ActiveWorkbook
$wsLogin = $wb.Worksheets.Item("LoginData")
$wsSearch = $wb.Worksheets.Item("SearchData")

# --- LoginData sheet: update key/value rows ---
# Row 2: osandaEmail / osanda@mailinator.com
$wsLogin.Range("A2").Value = "osandaEmail"
$wsLogin.Range("B2").Value = "osanda@mailinator.com"
# Row 3: osandaPassword / 1qaz2wsx@
$wsLogin.Range("A3").Value = "osandaPassword"
$wsLogin.Range("B3").Value = "1qaz2wsx@"
# Row 4: osandaProfileName / Osanda Nimalarathna
$wsLogin.Range("A4").Value = "osandaProfileName"
$wsLogin.Range("B4").Value = "Osanda Nimalarathna"

# --- SearchData sheet: replace question/answer rows with dress search data ---
# Row 2: tShirtDress / Faded Short Sleeve T-shirts
$wsSearch.Range("A2").Value = "tShirtDress"
$wsSearch.Range("B2").Value = "Faded Short Sleeve T-shirts"
# Row 3: clear the now-redundant Password/welcome1! test data
[void]$wsSearch.Range("A3").ClearContents()
[void]$wsSearch.Range("B3").ClearContents()

# --- Selection / active tab updates ---
[void]$wsLogin.Range("I14").Select()
[void]$wsSearch.Range("F6").Select()
